# Update "想去人数" (F column) counts across sheets to match the
# regenerated data snapshot (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

function Set-FValue {
    param(
        [string]$SheetName,
        [int]$Row,
        [double]$NewValue
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Cells.Item($Row, 6).Value = $NewValue
}

# 展览 (Exhibitions)
Set-FValue "展览" 3 994
Set-FValue "展览" 4 248
Set-FValue "展览" 5 445
Set-FValue "展览" 6 699
Set-FValue "展览" 7 247
Set-FValue "展览" 9 25
Set-FValue "展览" 10 392
Set-FValue "展览" 11 197
Set-FValue "展览" 12 68
Set-FValue "展览" 13 800
Set-FValue "展览" 15 1968
Set-FValue "展览" 16 464
Set-FValue "展览" 17 6831
Set-FValue "展览" 18 510
Set-FValue "展览" 19 517
Set-FValue "展览" 21 87
Set-FValue "展览" 24 135

# 演出 (Performances)
Set-FValue "演出" 4 30
Set-FValue "演出" 13 52
Set-FValue "演出" 20 25

# 本地生活 (Local Life)
Set-FValue "本地生活" 2 5454
Set-FValue "本地生活" 4 378

# 全部类型 (All Types) - aggregated sheet
Set-FValue "全部类型" 3 5454
Set-FValue "全部类型" 5 378
Set-FValue "全部类型" 8 30
Set-FValue "全部类型" 11 994
Set-FValue "全部类型" 14 248
Set-FValue "全部类型" 15 445
Set-FValue "全部类型" 16 699
Set-FValue "全部类型" 17 247
Set-FValue "全部类型" 20 25
Set-FValue "全部类型" 21 392
Set-FValue "全部类型" 22 197
Set-FValue "全部类型" 24 68
Set-FValue "全部类型" 26 800
Set-FValue "全部类型" 29 1968
Set-FValue "全部类型" 30 464
Set-FValue "全部类型" 31 6831
Set-FValue "全部类型" 32 52
Set-FValue "全部类型" 33 510
Set-FValue "全部类型" 34 517
Set-FValue "全部类型" 36 87
Set-FValue "全部类型" 41 135
Set-FValue "全部类型" 46 25
